$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(7, 1).Value = 9407
$ws.Cells.Item(7, 2).Value = 9335.1200000000008
$ws.Cells.Item(7, 3).Value = 107.96
$ws.Cells.Item(7, 4).Value = 108.79
$ws.Cells.Item(7, 5).Value = $false
$ws.Cells.Item(7, 6).Value = 0.77
$ws.Cells.Item(7, 7).Value = 42609.48846064815
$ws.Cells.Item(7, 7).NumberFormat = "m/d/yy h:mm"
$ws.Cells.Item(7, 8).Value = $true
